# Historyboard_RH.pptx - "se agrego el archivo reportes y nuevos estilos"
#
# New styles: 5 of the full-width "row" rectangles on slide 1 get an
# explicit colored outline (red for the first three rows, green for the
# last two) so the rows stand out in the storyboard.

$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

# RGB() is not available in this host, so colors are passed as the
# packed Long PowerPoint/VBA uses internally: R + G*256 + B*65536.
$red   = 255      # RRGGBB = FF0000
$green = 5287936  # RRGGBB = 00B050

$redShapes = @("Rectángulo 8", "Rectángulo 20", "Rectángulo 34")
foreach ($name in $redShapes) {
    $shp = $s1.Shapes.Item($name)
    $shp.Line.Visible = $true
    $shp.Line.ForeColor.RGB = $red
}

$greenShapes = @("Rectángulo 47", "Rectángulo 58")
foreach ($name in $greenShapes) {
    $shp = $s1.Shapes.Item($name)
    $shp.Line.Visible = $true
    $shp.Line.ForeColor.RGB = $green
}
